# PatternRecognition-Schedule.xlsx update
#
# The "Progress legend" block (Progress / Not Started / Minor questions /
# Major review needed / Not noted / Done, with K = threshold, L = label)
# previously lived at K20:L25, right below the morning-session table. It is
# moved up to K3:L8, alongside the first rows of the schedule, and the old
# location is cleared out. The active-cell selection is also updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the legend block: K20:L25 -> K3:L8 (values + styles travel together).
$legendSource = $ws.Range("K20:L25")
$legendTarget = $ws.Range("K3:L8")
$legendSource.Copy($legendTarget)

# Remove the now-empty trailing legend cells (J/K/L) from rows 19-26 entirely
# (not just blank them out) so the row shrinks back to column I.
$ws.Range("J19:L26").Clear()

# Update the saved selection to reflect where the user ended up (C25).
$ws.Range("C25").Select() | Out-Null
